# Slide 1 ("Java语言程序设计" title slide): the college/school text box
# (shape id=21, "学院：...") is updated to the new college name, and the
# (auto-fit) text box is re-centered/shrunk to match the shorter text,
# keeping the same vertical position/height and the same horizontal
# center point as before.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)

# Left/Width expressed in points (EMU / 12700) chosen so they round-trip
# to the exact target EMU values (5312187 / 1783122) through the COM
# layer's single-precision Left/Width storage.
$sh.Left = 418.282470703125
$sh.Width = 140.4033203125

$sh.TextFrame.TextRange.Text = "学院：人工智能学院"
